# Defect Density.xlsx -- replace "Project 3" (Apache Bigtop) with
# "Apache Commons Configuration" and "Project 4" (Apache Cayenne) with
# "JFreeChart", including their per-version defect-density rows.
#
# NOTE: the order in which new text values are first written controls the
# order they land in the rebuilt shared-string table on save, so the text
# cells are written in the same order the original workbook's table ended
# up in: Project-3 name, Project-4 name, then the Project-4 version labels
# from newest (row 24) to oldest (row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B14").Value = "Apache Commons Configuration"
$ws.Range("B20").Value = "JFreeChart"

$ws.Range("C24").Value = "1.5.0"
$ws.Range("C23").Value = "1.0.19"
$ws.Range("C22").Value = "1.0.18"
$ws.Range("C21").Value = "1.0.17"
$ws.Range("C20").Value = "1.0.16"

# ---------------------------------------------------------------
# Project 3 block (rows 14-18): Apache Bigtop -> Apache Commons Configuration
# Version numbers switch from text labels to plain numbers.
# ---------------------------------------------------------------
$ws.Range("C14").Value = 2.3
$ws.Range("D14").Value = 13
$ws.Range("E14").Value = 88468
$ws.Range("F14").Value = 88.468

$ws.Range("C15").Value = 2.4
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = 88851
$ws.Range("F15").Value = 88.851

$ws.Range("C16").Value = 2.5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 89035
$ws.Range("F16").Value = 89.035

$ws.Range("C17").Value = 2.6
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 89269
$ws.Range("F17").Value = 89.269

$ws.Range("C18").Value = 2.7
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 88749
$ws.Range("F18").Value = 88.749

# ---------------------------------------------------------------
# Project 4 block (rows 20-24): Apache Cayenne -> JFreeChart
# Version numbers stay as text, and a 5th version row (24) is added.
# ---------------------------------------------------------------
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 308748
$ws.Range("F20").Value = 308.748

$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 311144
$ws.Range("F21").Value = 311.144

$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 318874
$ws.Range("F22").Value = 318.874

$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 319056
$ws.Range("F23").Value = 319.056

$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 296835
$ws.Range("F24").Value = 296.835
$ws.Range("G24").Formula = "=D24/F24"

# ---------------------------------------------------------------
# Column B needs to be a bit wider to fit "Apache Commons Configuration".
# ---------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 29
